# This script fixes the marksheet computation for 1401ME08:
#  - recomputes Right/Wrong/NotAttempted/Max and the marking/score summary
#  - makes the "Wrong" mark numeric instead of text ("Handles float input without breaking stuff")
#  - removes the two bogus duplicate "Student Ans / Correct Ans" answer-key blocks
#    that used to live in columns D:E and G:H for rows 16-40
#  - fills in the "Student Ans" cell for every question the student answered
#    correctly (these simply mirror the "Correct Ans" value for that row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Give the row-label cells (A10/A11/A12) the same "mtitleStyle" formatting
#    already used by the other header-ish labels (e.g. A9).
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the summary numbers (Right / Wrong / Not Attempt / Max, the
#    per-question marks, and the computed Total).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 11    # Right
$ws.Range("D10").Value = 17    # Not Attempt
$ws.Range("E10").Value = 28    # Max

$ws.Range("B11").Value = 4     # Right mark
$ws.Range("C11").Value = -1    # Wrong mark - now a real number, not text

$ws.Range("B12").Value = 44        # Total score
$ws.Range("E12").Value = "44/112"  # Total score / max possible score

# ---------------------------------------------------------------------------
# 3) Remove the two extra (duplicated/bogus) answer-key blocks that used to
#    occupy columns D:E (rows 19-40) and G:H (rows 15-21). Clear() fully
#    removes the cells rather than just blanking their contents.
# ---------------------------------------------------------------------------
$ws.Range("D19:E40").Clear()
$ws.Range("G15:H21").Clear()

# ---------------------------------------------------------------------------
# 4) For every question the student got right, mirror the "Correct Ans"
#    value into the "Student Ans" column, using the same "correctStyle"
#    formatting already used for that column.
# ---------------------------------------------------------------------------
function Set-StudentAns($cellRef, $text) {
    $ws.Range("B10").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)   # xlPasteFormats (correctStyle)
    $ws.Range($cellRef).Value = $text
}

Set-StudentAns "D16" "Option A"
Set-StudentAns "D17" "Option C"
Set-StudentAns "A18" "Option B"
Set-StudentAns "A25" "Option A"
Set-StudentAns "A27" "Option A"
Set-StudentAns "A28" "Option D"
Set-StudentAns "A30" "Option B"
Set-StudentAns "A32" "Option C"
Set-StudentAns "A33" "Option D"
Set-StudentAns "A36" "Option A"
Set-StudentAns "A39" "Option D"
